# update test and example files (APPL should always have been AAPL)
#
# The placeholder single-letter symbols "A" and "B" used in the example
# workbook are replaced by the more descriptive placeholder symbols
# "ALOTOFLOVE" and "BERIGHTBACK" respectively.

$wb = $excel.ActiveWorkbook

$buyOrders  = $wb.Worksheets.Item("buy_orders")
$sellOrders = $wb.Worksheets.Item("sell_orders")

# buy_orders: B2 "A" -> "ALOTOFLOVE", B3 "B" -> "BERIGHTBACK"
$buyOrders.Range("B2").Value = "ALOTOFLOVE"
$buyOrders.Range("B3").Value = "BERIGHTBACK"

# sell_orders: B2 "A" -> "ALOTOFLOVE"
$sellOrders.Range("B2").Value = "ALOTOFLOVE"

# Restore the selection/active sheet back to the first sheet, and leave
# a selection on B2 on both edited sheets, matching the saved view state.
$sellOrders.Activate()
$sellOrders.Range("B2").Select() | Out-Null
$buyOrders.Activate()
$buyOrders.Range("B2").Select() | Out-Null

$wb.Save()
